# Moved 0xdeadbeef to column D. (actually: DEADBEEF.1 / Out-of-scope markers
# move from column D to a new column E on the "Steps Overview" sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps Overview")

# Give the new column E a width, matching the source column's sizing intent.
$ws.Columns.Item(5).ColumnWidth = 18.25

# Rows whose column-D cell holds a "DEADBEEF.n" / "Out of scope" marker that
# needs to move one column to the right, into the (new) column E.
$rows = @(116,117,118,119,120,121,122,123,124,125,127,128,129,130,131,132,211,222,224,225,227,228,329,334,340,356,362,378,384,414,420,422,423,424,426,427,437,443,449)

foreach ($r in $rows) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value2 = $dCell.Value2
    $dCell.ClearContents()
}

# Restore the selection/active cell the author left the sheet in.
$ws.Range("E422:E427").Select()
$ws.Application.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null

# The "SP800-78" sheet picked up a fresh A1:E4 selection too.
$ws2 = $wb.Worksheets.Item("SP800-78")
$ws2.Range("A1:E4").Select()

# Re-select the Steps Overview sheet / cell the author ended on.
$ws.Select()
$ws.Range("E422:E427").Select()
